$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 3.897971727105102
$ws.Range("B3").Value = 2.263060848186847
$ws.Range("C3").Value = 2.338660110633067
$ws.Range("D3").Value = 0.7713847572218808
$ws.Range("E3").Value = 68.44253226797787
$ws.Range("F3").Value = 0.2904585125998771
